$d = $word.ActiveDocument

# Locate the anchor paragraph that ends the bibliography section
# ("Thomson Pioneira (2008).") and the paragraph that holds the trailing
# copyright/footer line ("... Creative Commons Attribution"). Everything
# in between (an empty paragraph, the "Ver no Jupiter ..." paragraph, and
# the copyright paragraph itself) must be removed as whole paragraphs,
# leaving the blank paragraph that precedes the page break untouched.
$anchorIdx = 0
$footerEndIdx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Thomson Pioneira (2008).")) {
        $anchorIdx = $i
    }
    if ($p.Range.Text.Contains("Creative Commons Attribution")) {
        $footerEndIdx = $i
    }
    $i = $i + 1
}

if ($anchorIdx -gt 0 -and $footerEndIdx -gt $anchorIdx) {
    $startIdx = $anchorIdx + 1
    $deleteStart = $d.Paragraphs.Item($startIdx).Range.Start
    $deleteEnd = $d.Paragraphs.Item($footerEndIdx + 1).Range.Start
    $r = $d.Range($deleteStart, $deleteEnd)
    $r.Delete()
}
